# Re-update: Games (Zixel Technologies)
# Slide 1 has two "breadcrumb" textboxes that describe a URL path.
# TextBox 5 (shape 1): "…/palette"          -> "…/tic-tac-toe"
# TextBox 6 (shape 2): "…/play/tic-tac-toe" -> "…/"
# Each textbox auto-fits its width to the text (a:spAutoFit / wrap="none"),
# so the stored width (EMU) must follow the new text length.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$shape1 = $s.Shapes.Item(1)
$shape1.TextFrame.TextRange.Text = "…/tic-tac-toe"
$shape1.Width = 2710999 / 12700

$shape2 = $s.Shapes.Item(2)
$shape2.TextFrame.TextRange.Text = "…/"
$shape2.Width = 729687 / 12700
